$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated simulation values (dates shifted by one day; rows 9-50 recomputed
# charging power / SOC / net energy after adding sample script + plot ylim changes).
$ws.Cells.Item(2, 1).Value = 45689.375
$ws.Cells.Item(3, 1).Value = 45689.38541666666
$ws.Cells.Item(4, 1).Value = 45689.39583333334
$ws.Cells.Item(5, 1).Value = 45689.40625
$ws.Cells.Item(6, 1).Value = 45689.41666666666
$ws.Cells.Item(7, 1).Value = 45689.42708333334
$ws.Cells.Item(8, 1).Value = 45689.4375
$ws.Cells.Item(9, 1).Value = 45689.44791666666 ; $ws.Cells.Item(9, 2).Value = 11 ; $ws.Cells.Item(9, 3).Value = 13.66666666666667 ; $ws.Cells.Item(9, 4).Value = 1.833333333333333
$ws.Cells.Item(10, 1).Value = 45689.45833333334 ; $ws.Cells.Item(10, 3).Value = 18.80000000000001 ; $ws.Cells.Item(10, 4).Value = 4.399999999999999
$ws.Cells.Item(11, 1).Value = 45689.46875 ; $ws.Cells.Item(11, 3).Value = 18.80000000000001 ; $ws.Cells.Item(11, 4).Value = 4.399999999999999
$ws.Cells.Item(12, 1).Value = 45689.47916666666 ; $ws.Cells.Item(12, 3).Value = 18.80000000000001 ; $ws.Cells.Item(12, 4).Value = 4.399999999999999
$ws.Cells.Item(13, 1).Value = 45689.48958333334 ; $ws.Cells.Item(13, 3).Value = 18.80000000000001 ; $ws.Cells.Item(13, 4).Value = 4.399999999999999
$ws.Cells.Item(14, 1).Value = 45689.5 ; $ws.Cells.Item(14, 3).Value = 18.80000000000001 ; $ws.Cells.Item(14, 4).Value = 4.399999999999999
$ws.Cells.Item(15, 1).Value = 45689.51041666666 ; $ws.Cells.Item(15, 3).Value = 18.80000000000001 ; $ws.Cells.Item(15, 4).Value = 4.399999999999999
$ws.Cells.Item(16, 1).Value = 45689.52083333334 ; $ws.Cells.Item(16, 3).Value = 18.80000000000001 ; $ws.Cells.Item(16, 4).Value = 4.399999999999999
$ws.Cells.Item(17, 1).Value = 45689.53125 ; $ws.Cells.Item(17, 3).Value = 18.80000000000001 ; $ws.Cells.Item(17, 4).Value = 4.399999999999999
$ws.Cells.Item(18, 1).Value = 45689.54166666666 ; $ws.Cells.Item(18, 3).Value = 18.80000000000001 ; $ws.Cells.Item(18, 4).Value = 4.399999999999999
$ws.Cells.Item(19, 1).Value = 45689.55208333334 ; $ws.Cells.Item(19, 3).Value = 18.80000000000001 ; $ws.Cells.Item(19, 4).Value = 4.399999999999999
$ws.Cells.Item(20, 1).Value = 45689.5625 ; $ws.Cells.Item(20, 3).Value = 18.80000000000001 ; $ws.Cells.Item(20, 4).Value = 4.399999999999999
$ws.Cells.Item(21, 1).Value = 45689.57291666666 ; $ws.Cells.Item(21, 3).Value = 18.80000000000001 ; $ws.Cells.Item(21, 4).Value = 4.399999999999999
$ws.Cells.Item(22, 1).Value = 45689.58333333334 ; $ws.Cells.Item(22, 3).Value = 18.80000000000001 ; $ws.Cells.Item(22, 4).Value = 4.399999999999999
$ws.Cells.Item(23, 1).Value = 45689.59375 ; $ws.Cells.Item(23, 3).Value = 18.80000000000001 ; $ws.Cells.Item(23, 4).Value = 4.399999999999999
$ws.Cells.Item(24, 1).Value = 45689.60416666666 ; $ws.Cells.Item(24, 3).Value = 18.80000000000001 ; $ws.Cells.Item(24, 4).Value = 4.399999999999999
$ws.Cells.Item(25, 1).Value = 45689.61458333334 ; $ws.Cells.Item(25, 3).Value = 18.80000000000001 ; $ws.Cells.Item(25, 4).Value = 4.399999999999999
$ws.Cells.Item(26, 1).Value = 45689.625 ; $ws.Cells.Item(26, 3).Value = 18.80000000000001 ; $ws.Cells.Item(26, 4).Value = 4.399999999999999
$ws.Cells.Item(27, 1).Value = 45689.63541666666 ; $ws.Cells.Item(27, 3).Value = 18.80000000000001 ; $ws.Cells.Item(27, 4).Value = 4.399999999999999
$ws.Cells.Item(28, 1).Value = 45689.64583333334 ; $ws.Cells.Item(28, 3).Value = 18.80000000000001 ; $ws.Cells.Item(28, 4).Value = 4.399999999999999
$ws.Cells.Item(29, 1).Value = 45689.65625 ; $ws.Cells.Item(29, 3).Value = 18.80000000000001 ; $ws.Cells.Item(29, 4).Value = 4.399999999999999
$ws.Cells.Item(30, 1).Value = 45689.66666666666 ; $ws.Cells.Item(30, 3).Value = 18.80000000000001 ; $ws.Cells.Item(30, 4).Value = 4.399999999999999
$ws.Cells.Item(31, 1).Value = 45689.67708333334 ; $ws.Cells.Item(31, 3).Value = 18.80000000000001 ; $ws.Cells.Item(31, 4).Value = 4.399999999999999
$ws.Cells.Item(32, 1).Value = 45689.6875 ; $ws.Cells.Item(32, 3).Value = 18.80000000000001 ; $ws.Cells.Item(32, 4).Value = 4.399999999999999
$ws.Cells.Item(33, 1).Value = 45689.69791666666 ; $ws.Cells.Item(33, 3).Value = 18.80000000000001 ; $ws.Cells.Item(33, 4).Value = 4.399999999999999
$ws.Cells.Item(34, 1).Value = 45689.70833333334 ; $ws.Cells.Item(34, 3).Value = 18.80000000000001 ; $ws.Cells.Item(34, 4).Value = 4.399999999999999
$ws.Cells.Item(35, 1).Value = 45689.71875 ; $ws.Cells.Item(35, 3).Value = 18.80000000000001 ; $ws.Cells.Item(35, 4).Value = 4.399999999999999
$ws.Cells.Item(36, 1).Value = 45689.72916666666 ; $ws.Cells.Item(36, 3).Value = 18.80000000000001 ; $ws.Cells.Item(36, 4).Value = 4.399999999999999
$ws.Cells.Item(37, 1).Value = 45689.73958333334 ; $ws.Cells.Item(37, 3).Value = 18.80000000000001 ; $ws.Cells.Item(37, 4).Value = 4.399999999999999
$ws.Cells.Item(38, 1).Value = 45689.75 ; $ws.Cells.Item(38, 3).Value = 18.80000000000001 ; $ws.Cells.Item(38, 4).Value = 4.399999999999999
$ws.Cells.Item(39, 1).Value = 45689.76041666666 ; $ws.Cells.Item(39, 3).Value = 18.80000000000001 ; $ws.Cells.Item(39, 4).Value = 4.399999999999999
$ws.Cells.Item(40, 1).Value = 45689.77083333334 ; $ws.Cells.Item(40, 3).Value = 18.80000000000001 ; $ws.Cells.Item(40, 4).Value = 4.399999999999999
$ws.Cells.Item(41, 1).Value = 45689.78125 ; $ws.Cells.Item(41, 3).Value = 18.80000000000001 ; $ws.Cells.Item(41, 4).Value = 4.399999999999999
$ws.Cells.Item(42, 1).Value = 45689.79166666666 ; $ws.Cells.Item(42, 3).Value = 18.80000000000001 ; $ws.Cells.Item(42, 4).Value = 4.399999999999999
$ws.Cells.Item(43, 1).Value = 45689.80208333334 ; $ws.Cells.Item(43, 3).Value = 18.80000000000001 ; $ws.Cells.Item(43, 4).Value = 4.399999999999999
$ws.Cells.Item(44, 1).Value = 45689.8125 ; $ws.Cells.Item(44, 3).Value = 18.80000000000001 ; $ws.Cells.Item(44, 4).Value = 4.399999999999999
$ws.Cells.Item(45, 1).Value = 45689.82291666666 ; $ws.Cells.Item(45, 3).Value = 18.80000000000001 ; $ws.Cells.Item(45, 4).Value = 4.399999999999999
$ws.Cells.Item(46, 1).Value = 45689.83333333334 ; $ws.Cells.Item(46, 3).Value = 18.80000000000001 ; $ws.Cells.Item(46, 4).Value = 4.399999999999999
$ws.Cells.Item(47, 1).Value = 45689.84375 ; $ws.Cells.Item(47, 3).Value = 18.80000000000001 ; $ws.Cells.Item(47, 4).Value = 4.399999999999999
$ws.Cells.Item(48, 1).Value = 45689.85416666666 ; $ws.Cells.Item(48, 3).Value = 18.80000000000001 ; $ws.Cells.Item(48, 4).Value = 4.399999999999999
$ws.Cells.Item(49, 1).Value = 45689.86458333334 ; $ws.Cells.Item(49, 3).Value = 18.80000000000001 ; $ws.Cells.Item(49, 4).Value = 4.399999999999999
$ws.Cells.Item(50, 1).Value = 45689.875 ; $ws.Cells.Item(50, 3).Value = 18.80000000000001 ; $ws.Cells.Item(50, 4).Value = 4.399999999999999
